$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source Information")

# Step 1: remove the "Sales vSalesPerson" row (old row 11); this shifts rows 12-14 up to 11-13
# and the table (ListObject) auto-shrinks its range to A1:I13.
$ws.Rows.Item(11).Delete()

# Step 2: renumber the "Table No" column (A) for the remaining data rows.
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 6
$ws.Range("A5").Value = 8
$ws.Range("A6").Value = 10
$ws.Range("A7").Value = 12
$ws.Range("A8").Value = 14
$ws.Range("A9").Value = 16
$ws.Range("A10").Value = 18
$ws.Range("A11").Value = 21
$ws.Range("A12").Value = 23
$ws.Range("A13").Value = 25

# Step 3: update the "Modification Description" (column I) wording for rows 2,3,5-10.
$tmp = @'
1. This changes the data types of the columns in the source table to "Table No" as Int64, "Table Name" as text, "Table Mode" as text, "Table Type" as text, "Table Source" as text, "Original Table Name" as text, "Table Query" as




'@
$ws.Range("I2").Value = $tmp

$tmp = @'
1. Change the data types of the columns in Table1 to Text for "Measure Name", "Measure Expression", "Measure Data Type", and "Measure Description".




'@
$ws.Range("I3").Value = $tmp

$tmp = @'
1. This changes the data type of the column "No relationships present in this file" to "any".




'@
$ws.Range("I5").Value = $tmp

$tmp = @'
1. This sentence changes the data types of each column in the promoted headers table to different types, such as Int64 and text.




'@
$ws.Range("I6").Value = $tmp

$tmp = @'
1. This sentence changes the data types of several columns in the source table to Int64, text, and text respectively.




'@
$ws.Range("I7").Value = $tmp

$tmp = @'
1. This changes the data type of the column "No measures presented in this file" to any type.




'@
$ws.Range("I8").Value = $tmp

$tmp = @'
1. The command is changing each column's type in a table to a specified type.




'@
$ws.Range("I9").Value = $tmp

$tmp = @'
1. This changes the data type of the column labeled "No relationships present in this file" to any.




'@
$ws.Range("I10").Value = $tmp

# Step 4: the "data" row (now row 13, previously row 14) keeps its H13 text, but its
# Modification Description (I13) is reworded.
$tmp = @'
1. Promoted Headers is a function that takes the data from the provided source, and creates an Excel table with all scalar values as headers.



2. The column types in the table were changed to text, datetime, number and text for the respective columns.




'@
$ws.Range("I13").Value = $tmp

# Step 5: grow the "Source" table reference to A1:I26 (and its AutoFilter range to match),
# as in the target workbook -- the extra rows stay empty (no cell data), matching the
# worksheet dimension staying at A1:I13.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I26"))

